{"js": "// The document starts as 11 paragraphs:\n//   1. \"kubectl\"\n//   2. \"\" (empty)\n//   3. \"kubectl config get-contexts\"\n//   4. \"\" (empty)\n//   5. \"kubectl create deployment nginx --image=nginx\"\n//   6. \"kubectl expose deployment nginx --type=NodePort --port=80\"\n//   7. \"kubectl get service nginx\"\n//   8. \"\" (empty)\n//   9. \"\" (empty)\n//  10. \"\" (empty)\n//  11. \"\" (empty)\n//\n// Target (per the commit's diff):\n//   1. \"Ol\u00e1!\"\n//   2. \"Meu nome \u00e9 #nome# \"\n//   3. \"E eu tenho #idade# anos\"   (this is the old paragraph 11, with\n//      every paragraph between it and paragraph 2 removed)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Replace the text of the first paragraph: \"kubectl\" -> \"Ol\u00e1!\"\nitems[0].insertText(\"Ol\u00e1!\", Word.InsertLocation.replace);\n\n// The second paragraph was an empty run; give it the new greeting line.\nitems[1].insertText(\"Meu nome \u00e9 #nome# \", Word.InsertLocation.replace);\n\n// The very last paragraph survives (also previously empty) and receives\n// the final templated sentence.\nconst lastIndex = items.length - 1;\nitems[lastIndex].insertText(\"E eu tenho #idade# anos\", Word.InsertLocation.replace);\n\n// Every paragraph strictly between the 2nd and the last one is removed\n// (the \"kubectl ...\" command lines plus the blank spacer paragraphs).\nfor (let i = lastIndex - 1; i >= 2; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The document starts as 11 paragraphs:\n#   1. \"kubectl\"\n#   2. \"\" (empty)\n#   3. \"kubectl config get-contexts\"\n#   4. \"\" (empty)\n#   5. \"kubectl create deployment nginx --image=nginx\"\n#   6. \"kubectl expose deployment nginx --type=NodePort --port=80\"\n#   7. \"kubectl get service nginx\"\n#   8. \"\" (empty)\n#   9. \"\" (empty)\n#  10. \"\" (empty)\n#  11. \"\" (empty)\n#\n# Target (per the commit's diff):\n#   1. \"Ol\u00e1!\"\n#   2. \"Meu nome \u00e9 #nome# \"\n#   3. \"E eu tenho #idade# anos\"   (this is the old paragraph 11, with\n#      every paragraph between it and paragraph 2 removed)\n\n$d = $word.ActiveDocument\n\n# Replace the text of the first paragraph: \"kubectl\" -> \"Ol\u00e1!\"\n$d.Paragraphs.Item(1).Range.Text = \"Ol\u00e1!\"\n\n# The second paragraph was an empty run; give it the new greeting line.\n$d.Paragraphs.Item(2).Range.Text = \"Meu nome \u00e9 #nome# \"\n\n# Remove every paragraph strictly between the 2nd and the last one\n# (the \"kubectl ...\" command lines plus the blank spacer paragraphs).\n# Walk backwards so earlier indices stay valid as later ones are removed.\n$lastIndex = $d.Paragraphs.Count\nfor ($i = $lastIndex - 1; $i -ge 3; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# The very last paragraph survives (also previously empty) and receives\n# the final templated sentence.\n$d.Paragraphs.Item(3).Range.Text = \"E eu tenho #idade# anos\"\n"}
